$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - "Subtitle 2" placeholder: midterm-presentation feedback update
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$tr1 = $subtitle.TextFrame.TextRange

# Paragraph 2 : "27th November 2017" -> "28th November 2017"
# (re-seed with a throwaway value first so the whole run is rewritten as a
#  single run instead of only the differing characters being patched)
$para2 = $tr1.Paragraphs(2)
$para2.Text = "Z"
$para2.Text = "28th November 2017"

# Paragraph 3 : merge "Max " / "Premi" / ", Master Semester Project (12 credits)"
# into a single run reading "Max Premi, Master Semester Project (12 credits)"
$para3 = $tr1.Paragraphs(3)
$para3.Text = "Z"
$para3.Text = "Max Premi, Master Semester Project (12 credits)"

# Paragraph 4 : move the hyphen from ", Juan Ramon Troncoso-" onto "Pastoriza"
# so the runs become ", Juan Ramon " / "Troncoso-Pastoriza"
$para4 = $tr1.Paragraphs(4)
$full4 = $para4.Text
$idxComma = $full4.IndexOf(", Juan Ramon Troncoso-")
$run4c = $para4.Characters($idxComma + 1, ", Juan Ramon Troncoso-".Length)
$run4c.Text = ", Juan Ramon "

$full4 = $para4.Text
$idxPastoriza = $full4.IndexOf("Pastoriza")
$run4d = $para4.Characters($idxPastoriza + 1, "Pastoriza".Length)
$run4d.Text = "Troncoso-Pastoriza"

# ---------------------------------------------------------------------------
# Slide 2 - content placeholder: project title wording
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$content2 = $s2.Shapes.Item(2)
$para2b = $content2.TextFrame.TextRange.Paragraphs(1)

$full = $para2b.Text
$idxPP = $full.IndexOf("Privacy-preserving")
$runPP = $para2b.Characters($idxPP + 1, 18)
$runPP.Text = "Privacy-Preserving"

$full = $para2b.Text
$idxSys = $full.IndexOf(" data sharing system")
$runSys = $para2b.Characters($idxSys + 15, 6)
$runSys.Text = "systems"

# ---------------------------------------------------------------------------
# Slide 3 - content placeholder: "with data in the Data Providers"
#           -> "with data from the Data Providers"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$content3 = $s3.Shapes.Item(2)
$para3b = $content3.TextFrame.TextRange.Paragraphs(3)

$full = $para3b.Text
$idxIn = $full.IndexOf(" data in the Data Providers")
$runIn = $para3b.Characters($idxIn + 7, 2)
$runIn.Text = "from"

# ---------------------------------------------------------------------------
# Slide 4 - content placeholder: "Proof of correctness using by SNIPs proof
#           and Proof for input validation"
#           -> "Proof of correctness and Input validation  using SNIPs proof"
#              (with "SNIP" underlined and split from the trailing "s")
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$content4 = $s4.Shapes.Item(2)
$para4b = $content4.TextFrame.TextRange.Paragraphs(9)

# " by " -> " "
$full = $para4b.Text
$idxBy = $full.IndexOf(" by ")
$runBy = $para4b.Characters($idxBy + 1, 4)
$runBy.Text = " "

# "SNIPs" -> "SNIP" (underlined) + "s"
$full = $para4b.Text
$idxSnip = $full.IndexOf("SNIPs")
$runSnip = $para4b.Characters($idxSnip + 1, 4)
$runSnip.Font.Underline = $true

# " proof and Proof for input validation" -> " proof"
$full = $para4b.Text
$idxTail = $full.IndexOf(" proof and Proof for input validation")
$runTail = $para4b.Characters($idxTail + 1, 38)
$runTail.Text = " proof"

# " " (between "correctness" and "using") -> " and Input validation  "
$full = $para4b.Text
$idxMid = $full.IndexOf("correctness") + "correctness".Length
$runMid = $para4b.Characters($idxMid + 1, 1)
$runMid.Text = " and Input validation  "

# ---------------------------------------------------------------------------
# Slide 5 - content placeholder: "Server reconstruct from shares..."
#           -> "Servers reconstruct from shares..."
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$content5 = $s5.Shapes.Item(5)
$para5b = $content5.TextFrame.TextRange.Paragraphs(7)

$full5 = $para5b.Text
$idxServer = $full5.IndexOf("Server ")
$runServer = $para5b.Characters($idxServer + 1, "Server ".Length)
$runServer.Text = "Servers "
